# Apply the workbook update described by the commit "Update gh-pages to
# output generated at 456a3b4".
#
# Summary of the change:
#   1. Sheet "展览" (Exhibitions): several "want to go" counts (column F)
#      increased for a number of existing events.
#   2. Sheet "演出" (Performances): two events dated 2024-04-19
#      ("2024武侠新国风沉浸音乐会《射雕英雄传》" and
#      "动漫钢琴鬼才Kyle Xian互动演奏会（取消）") were removed, so every
#      later row shifts up by two rows.
#   3. Sheet "本地生活" (Local life): unchanged.
#   4. Sheet "全部类型" (All types, union of the other sheets): the same
#      two 2024-04-19 events were removed (rows shift up by two), and the
#      same "want to go" counts as in "展览" increased at their new row
#      positions.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Sheet "展览" - bump column F ("想去人数") values
# ---------------------------------------------------------------------
$wsExhibit = $wb.Worksheets.Item("展览")

$exhibitUpdates = @{
    3  = 1008
    4  = 26
    5  = 1177
    6  = 979
    7  = 296
    10 = 917
    11 = 334
    12 = 596
    13 = 536
    16 = 1298
    17 = 2955
    18 = 353
    19 = 1579
    20 = 1328
    21 = 767
    23 = 1314
    24 = 245
    26 = 1089
    27 = 382
    28 = 3368
    29 = 655
    31 = 1493
}

foreach ($row in $exhibitUpdates.Keys) {
    $wsExhibit.Cells.Item($row, 6).Value = $exhibitUpdates[$row]
}

# ---------------------------------------------------------------------
# 2. Sheet "演出" - remove the two 2024-04-19 rows (old rows 2 and 3)
# ---------------------------------------------------------------------
$wsShow = $wb.Worksheets.Item("演出")
$wsShow.Range("A2:I3").EntireRow.Delete()

# Renumber column A (the 0-based sequence index) for the remaining rows.
for ($r = 2; $r -le 15; $r++) {
    $wsShow.Cells.Item($r, 1).Value = $r - 1
}

# ---------------------------------------------------------------------
# 3. Sheet "本地生活" - no changes
# ---------------------------------------------------------------------

# ---------------------------------------------------------------------
# 4. Sheet "全部类型" - remove the same two 2024-04-19 rows (old rows 4
#    and 5), then bump the same "want to go" counts at their new
#    (shifted) row positions.
# ---------------------------------------------------------------------
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("A4:I5").EntireRow.Delete()

# Renumber column A for the remaining rows.
for ($r = 2; $r -le 47; $r++) {
    $wsAll.Cells.Item($r, 1).Value = $r - 1
}

$allUpdates = @{
    4  = 1008
    6  = 26
    7  = 1177
    8  = 979
    9  = 296
    21 = 917
    22 = 334
    23 = 596
    24 = 536
    27 = 1298
    28 = 2955
    29 = 353
    30 = 1579
    31 = 1328
    32 = 767
    34 = 1314
    35 = 245
    39 = 1089
    40 = 382
    41 = 3368
    42 = 655
    44 = 1493
}

foreach ($row in $allUpdates.Keys) {
    $wsAll.Cells.Item($row, 6).Value = $allUpdates[$row]
}
